$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..19, columns B,D,E,F,G,H,I
# (C column - Metodo - is unchanged: MP / MP_Aitken alternating)
$data = @{
    2  = @("fidap005.mtx", 4359440.548861713, 45, 0.001310110092163086, 27, "real", "general")
    3  = @("fidap005.mtx", 4359599.815813746, 30, 0.0007317066192626953, 27, "real", "general")
    4  = @("gr_30_30.mtx", 11.86393030934931, 113, 0.2811086177825928, 900, "integer", "symmetric")
    5  = @("gr_30_30.mtx", 11.86933820735485, 49, 0.129218578338623, 900, "integer", "symmetric")
    6  = @("05r0100.mtx", 12.41717072759173, 233, 0.01548933982849121, 236, "real", "general")
    7  = @("05r0100.mtx", 12.63404220515477, 37, 0.00245976448059082, 236, "real", "general")
    8  = @("bcsstm23.mtx", 8270300.23229669, 16, 0.07952570915222168, 3134, "real", "symmetric")
    9  = @("bcsstm23.mtx", 8270311.116852832, 11, 0.05410361289978027, 3134, "real", "symmetric")
    10 = @("dwt_1005.mtx", 16.40635975914216, 14, 0.006129026412963867, 1005, "pattern", "symmetric")
    11 = @("dwt_1005.mtx", 16.40643403711349, 12, 0.004851102828979492, 1005, "pattern", "symmetric")
    12 = @("herman4.mtx", 66.48767021742641, 100, 0.0426781177520752, 1104, "real", "general")
    13 = @("herman4.mtx", 66.5125181458945, 39, 0.01664400100708008, 1104, "real", "general")
    14 = @("dwt_1242.mtx", 9.38929622309684, 105, 0.0581505298614502, 1242, "pattern", "symmetric")
    15 = @("dwt_1242.mtx", 9.394188915458207, 37, 0.02043795585632324, 1242, "pattern", "symmetric")
    16 = @("dwt_1007.mtx", 8.912431480213071, 61, 0.0207819938659668, 1007, "pattern", "symmetric")
    17 = @("dwt_1007.mtx", 8.913562817502868, 45, 0.01305174827575684, 1007, "pattern", "symmetric")
    18 = @("dwt_245.mtx", 8.523646707394455, 48, 0.001848220825195312, 245, "pattern", "symmetric")
    19 = @("dwt_245.mtx", 8.523875352301671, 35, 0.001606464385986328, 245, "pattern", "symmetric")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - Matriz
    $ws.Cells.Item($row, 4).Value = $vals[1]   # D - Autovalor
    $ws.Cells.Item($row, 5).Value = $vals[2]   # E - Iteracoes
    $ws.Cells.Item($row, 6).Value = $vals[3]   # F - Tempo
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - Ordem
    $ws.Cells.Item($row, 8).Value = $vals[5]   # H - Campo
    $ws.Cells.Item($row, 9).Value = $vals[6]   # I - Simetria
}
